# Generate Report for Handoff
# Row 3 (the "b.md" row) moves from "Handed back" to "Ready for handoff":
# a new handoff package (b.63290e5768f688058c7b37413b0a5c26c308f864) was
# produced for both locales, so the Overview summary and each locale's
# detail sheet need their Status / Latest Handoff File / Latest Handoff
# Datetime columns refreshed for that row.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-17-09 10:17:14"

# ---- Helper: update a locale detail sheet's row 3 --------------------
function Update-LocaleSheet($sheetName, $handoffFile, $handoffDatetime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status
    $ws.Range("C3").Value = "Ready for handoff"

    # Latest Handoff File (cell text + underlying hyperlink display text;
    # the hyperlink target/address itself is left untouched)
    $ws.Range("D3").Value = $handoffFile
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Row -eq 3 -and $h.Range.Column -eq 4) {
            $h.TextToDisplay = $handoffFile
        }
    }

    # Latest Handoff Datetime
    $ws.Range("E3").Value = $handoffDatetime
}

Update-LocaleSheet "zh-cn" "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf" "2016-03-09 10:17:06"
Update-LocaleSheet "de-de" "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf" "2016-03-09 10:17:14"
